$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace e044 body (row 49) with the new "Sector" sub-event text,
# then insert three new rows for e044a / e044b / e044c right after it.

$e044Sector = @'
<Bold>e044 Panzerfaust Attack - Sector</Bold> 
<InlineUIContainer><Button Content='r15.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<InlineUIContainer><Button Content='Panzerfaust' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
Determine from which sector of the Battle Board attack is originating by rolling 1D according to 
<InlineUIContainer><Button Content='r5.12' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> :  
<InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21' > </Image></InlineUIContainer>
<LineBreak/><LineBreak/>
If the sector rolled is US controlled, no attack is made. If the sector is not US Controlled, a Panzerfaust marker is placed in the sector's close range.
<LineBreak/><LineBreak/>
'@

$e044aToAttack = @'
<Bold>e044a Panzerfaust Attack - To Attack</Bold> 
<InlineUIContainer><Button Content='r15.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<InlineUIContainer><Button Content='Panzerfaust' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
Determine if attack occurs based on scenario type, roll one die: 
<InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21' > </Image></InlineUIContainer>
<LineBreak/><LineBreak/>
  1 to 3 for Advance<LineBreak/>
  1 to 5 for Battle<LineBreak/>
  1 to 2 for Counterattack<LineBreak/>
<LineBreak/><LineBreak/>
<Underline>Modifiers:</Underline><LineBreak/>
  December 1944 or later = -1<LineBreak/>
  Sherman moving = -1<LineBreak/>
  Lead tank = -1<LineBreak/>
  Advancing Fire in zone = +3<LineBreak/>
  Attack from sector 1,2, or3 = -1<LineBreak/>
<LineBreak/><LineBreak/>
'@

$e044bToHit = @'
<Bold>e044b Panzerfaust Attack - To Hit</Bold> 
<InlineUIContainer><Button Content='r15.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<InlineUIContainer><Button Content='Panzerfaust' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<LineBreak/><LineBreak/>
Since an attack occurred, roll to hit:   
<InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21' > </Image></InlineUIContainer>
<LineBreak/><LineBreak/>
  1 to 7 is Hit<LineBreak/>
  8 to 10 is Miss<LineBreak/>
<LineBreak/><LineBreak/>
<Underline>Modifiers:</Underline><LineBreak/>
  Sherman moving = +2<LineBreak/>
  Advancing Fire in zone = +3<LineBreak/>
<LineBreak/><LineBreak/>
'@

$e044cToKill = @'
<Bold>e044c Panzerfaust Attack - To Kill</Bold> 
<InlineUIContainer><Button Content='r15.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   
<InlineUIContainer><Button Content='Panzerfaust' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
 <LineBreak/><LineBreak/>
Since there was a hit, roll again to see if your tank is knocked out (KO):  
<InlineUIContainer><Image Name='DieRollWhite' Height='21' Width='21'> </Image></InlineUIContainer>
<LineBreak/><LineBreak/>
  1 to 8  is KO<LineBreak/>
  9 to 10 is No Effect<LineBreak/>
<LineBreak/><LineBreak/>
'@

# Row 49 keeps label "e044" but its body becomes the Sector text.
$ws.Range("B49").Value = $e044Sector
$ws.Rows.Item(49).RowHeight = 150

# Insert three fresh rows below row 49 for e044a, e044b, e044c.
$ws.Rows.Item(50).Insert()
$ws.Rows.Item(51).Insert()
$ws.Rows.Item(52).Insert()

$ws.Range("A50").Value = "e044a"
$ws.Range("B50").Value = $e044aToAttack
$ws.Rows.Item(50).RowHeight = 270

$ws.Range("A51").Value = "e044b"
$ws.Range("B51").Value = $e044bToHit
$ws.Rows.Item(51).RowHeight = 210

$ws.Range("A52").Value = "e044c"
$ws.Range("B52").Value = $e044cToKill
$ws.Rows.Item(52).RowHeight = 150

# Restore the view to roughly where the author left it.
$ws.Range("B52").Select()

